# Update NATMI LR-pair TPM-derived values (Nid1-Itgb3) per new TPM re-computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Cells.Item(2, 7).Value = 28.743868
$ws.Cells.Item(2, 8).Value = 86.231604
$ws.Cells.Item(2, 9).Value = 0.0554303735704667
$ws.Cells.Item(2, 10).Value = 0.0554303735704667
$ws.Cells.Item(2, 13).Value = 0.2901893333333334
$ws.Cells.Item(2, 14).Value = 0.870568
$ws.Cells.Item(2, 15).Value = 0.03429389578125064
$ws.Cells.Item(2, 16).Value = 0.03429389578125064
$ws.Cells.Item(2, 17).Value = 8.341163892341335
$ws.Cells.Item(2, 18).Value = 75.07047503107201
$ws.Cells.Item(2, 19).Value = 0.001900923454341375
$ws.Cells.Item(2, 20).Value = 0.001900923454341375

# Row 3 (ECs -> FAPs)
$ws.Cells.Item(3, 7).Value = 28.743868
$ws.Cells.Item(3, 8).Value = 86.231604
$ws.Cells.Item(3, 9).Value = 0.0554303735704667
$ws.Cells.Item(3, 10).Value = 0.0554303735704667
$ws.Cells.Item(3, 15).Value = 0.8402845891331153
$ws.Cells.Item(3, 16).Value = 0.8402845891331153
$ws.Cells.Item(3, 17).Value = 204.3789809963787
$ws.Cells.Item(3, 18).Value = 1839.410828967408
$ws.Cells.Item(3, 19).Value = 0.0465772886811547
$ws.Cells.Item(3, 20).Value = 0.0465772886811547

# Row 4 (ECs -> MuSCs)
$ws.Cells.Item(4, 7).Value = 28.743868
$ws.Cells.Item(4, 8).Value = 86.231604
$ws.Cells.Item(4, 9).Value = 0.0554303735704667
$ws.Cells.Item(4, 10).Value = 0.0554303735704667
$ws.Cells.Item(4, 15).Value = 0.1254215150856341
$ws.Cells.Item(4, 16).Value = 0.1254215150856341
$ws.Cells.Item(4, 17).Value = 30.50576171421734
$ws.Cells.Item(4, 18).Value = 274.551855427956
$ws.Cells.Item(4, 19).Value = 0.006952161434970624
$ws.Cells.Item(4, 20).Value = 0.006952161434970623

# Row 5 (FAPs -> ECs)
$ws.Cells.Item(5, 9).Value = 0.848161237947095
$ws.Cells.Item(5, 10).Value = 0.8481612379470951
$ws.Cells.Item(5, 13).Value = 0.2901893333333334
$ws.Cells.Item(5, 14).Value = 0.870568
$ws.Cells.Item(5, 15).Value = 0.03429389578125064
$ws.Cells.Item(5, 16).Value = 0.03429389578125064
$ws.Cells.Item(5, 17).Value = 127.6313226331422
$ws.Cells.Item(5, 18).Value = 1148.68190369828
$ws.Cells.Item(5, 19).Value = 0.0290867530998542
$ws.Cells.Item(5, 20).Value = 0.0290867530998542

# Row 6 (FAPs -> FAPs)
$ws.Cells.Item(6, 9).Value = 0.848161237947095
$ws.Cells.Item(6, 10).Value = 0.8481612379470951
$ws.Cells.Item(6, 15).Value = 0.8402845891331153
$ws.Cells.Item(6, 16).Value = 0.8402845891331153
$ws.Cells.Item(6, 19).Value = 0.7126968173470092
$ws.Cells.Item(6, 20).Value = 0.7126968173470093

# Row 7 (FAPs -> MuSCs)
$ws.Cells.Item(7, 9).Value = 0.848161237947095
$ws.Cells.Item(7, 10).Value = 0.8481612379470951
$ws.Cells.Item(7, 15).Value = 0.1254215150856341
$ws.Cells.Item(7, 16).Value = 0.1254215150856341
$ws.Cells.Item(7, 18).Value = 4201.022410293065
$ws.Cells.Item(7, 19).Value = 0.1063776675002317
$ws.Cells.Item(7, 20).Value = 0.1063776675002317

# Row 8 (MuSCs -> ECs)
$ws.Cells.Item(8, 7).Value = 49.99334866666666
$ws.Cells.Item(8, 9).Value = 0.09640838848243828
$ws.Cells.Item(8, 10).Value = 0.09640838848243828
$ws.Cells.Item(8, 13).Value = 0.2901893333333334
$ws.Cells.Item(8, 14).Value = 0.870568
$ws.Cells.Item(8, 15).Value = 0.03429389578125064
$ws.Cells.Item(8, 16).Value = 0.03429389578125064
$ws.Cells.Item(8, 17).Value = 14.50753652068089
$ws.Cells.Item(8, 18).Value = 130.567828686128
$ws.Cells.Item(8, 19).Value = 0.003306219227055063
$ws.Cells.Item(8, 20).Value = 0.003306219227055063

# Row 9 (MuSCs -> FAPs)
$ws.Cells.Item(9, 7).Value = 49.99334866666666
$ws.Cells.Item(9, 9).Value = 0.09640838848243828
$ws.Cells.Item(9, 10).Value = 0.09640838848243828
$ws.Cells.Item(9, 15).Value = 0.8402845891331153
$ws.Cells.Item(9, 16).Value = 0.8402845891331153
$ws.Cells.Item(9, 17).Value = 355.4702400209324
$ws.Cells.Item(9, 19).Value = 0.08101048310495142
$ws.Cells.Item(9, 20).Value = 0.08101048310495142

# Row 10 (MuSCs -> MuSCs)
$ws.Cells.Item(10, 7).Value = 49.99334866666666
$ws.Cells.Item(10, 9).Value = 0.09640838848243828
$ws.Cells.Item(10, 10).Value = 0.09640838848243828
$ws.Cells.Item(10, 15).Value = 0.1254215150856341
$ws.Cells.Item(10, 16).Value = 0.1254215150856341
$ws.Cells.Item(10, 19).Value = 0.01209168615043181
$ws.Cells.Item(10, 20).Value = 0.0120916861504318
